$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ D = 0.8204873197414222; E = 0.654796511627907;  F = 0.5970841616964877; G = 0.6246100519930675 }
    3  = @{ D = 0.8183325045582629; E = 0.6495293265749457; F = 0.5944333996023857; G = 0.6207612456747404 }
    4  = @{ D = 0.8186640145864412; E = 0.6493880489560835; F = 0.5977468522200132; G = 0.6224982746721877 }
    5  = @{ D = 0.8211503397977788; E = 0.657856093979442;  F = 0.5937707090788602; G = 0.6241727621037966 }
    6  = @{ D = 0.8188297696005304; E = 0.6509433962264151; F = 0.5944333996023857; G = 0.6214063041219258 }
    7  = @{ D = 0.8206530747555113; E = 0.6554989075018208; F = 0.5964214711729622; G = 0.6245662734212353 }
    8  = @{ D = 0.8183325045582629; E = 0.6504005826656956; F = 0.5917826375082836; G = 0.61970853573907  }
    9  = @{ D = 0.8196585446709763; E = 0.6530909090909091; F = 0.5950960901259112; G = 0.6227461858529819 }
    10 = @{ D = 0.8196585446709763; E = 0.6517664023071377; F = 0.5990722332670643; G = 0.6243093922651933 }
    11 = @{ D = 0.8196585446709763; E = 0.6517664023071377; F = 0.5990722332670643; G = 0.6243093922651933 }
    12 = @{ D = 0.8199900546991546; E = 0.653372008701958;  F = 0.5970841616964877; G = 0.6239612188365651 }
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
